$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bus_Makhulu_f")
$ws.Range("F6").Formula = "=-0.0022+0.02"
$ws.Range("F6").NumberFormat = "0.00"
